$d = $word.ActiveDocument

# --- 1. Replace the figure caption text -----------------------------------
# Old:  "Nutrition radar (static with function group)."
# New:  "This figure shows the relative nutrient content of the six major
#        marine seafood groups. Only nutrient contents for nutrients with
#        EARs are shown."
$oldText = "Nutrition radar (static with function group)."
$newText = "This figure shows the relative nutrient content of the six major marine seafood groups. Only nutrient contents for nutrients with EARs are shown."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# --- 2. Re-seat the _GoBack bookmark in the middle of the new sentence ----
# In the source edit the last-edit marker ends up right after
# "...marine seafood g" (i.e. before "roups. Only nutrient ..."), which is
# 81 characters into the new sentence.
$bmOffset = 81
$paraStart = $d.Paragraphs(1).Range.Start
$bmRange = $d.Range($paraStart + $bmOffset, $paraStart + $bmOffset)
$d.Bookmarks.Add("_GoBack", $bmRange)
